$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "0.5a - two layered NN" entry is being inserted at row 29. The row that
# used to live at row 29 ("-", "-", "Create a tool to measure convergance
# speed", "In iterations") is pushed down to row 30, but its Notes cell is not
# carried over (left blank). Row 29 itself gets the brand-new entry.

# Old row 29 -> row 30 (Notes/D left empty). Pick up B30's format from the
# old B29 (date-style cell) before its value gets overwritten below.
$ws.Range("B29").Copy()
$ws.Range("B30").PasteSpecial(-4122)

$ws.Range("A30").Value = "-"
$ws.Range("B30").Value = "-"
$ws.Range("C30").Value = "Create a tool to measure convergance speed"
$ws.Range("D30").Value = ""

# New row 29 content: date 2016-01-24 (serial 42393), version 0.5a
$ws.Range("A29").Value = 42393
$ws.Range("B29").Value = "0.5a"
$ws.Range("C29").Value = "Create a two layered neural net"
$ws.Range("D29").Value = "Needed some help from example code for this, forgot to threshold the hidden layer" + [char]0x2026

# Reflect the final selection left by the edit
$ws.Range("D30").Select()
